# Adds a new "2022-Q4" sheet (with its fund-holding detail data) right after
# the "总计" summary sheet, pushing the former "2022-Q3" sheet (and every
# quarter after it) one position later, and inserts a new summary row for
# 2022-Q4 at the top of the data in the "总计" sheet.

$wb = $excel.ActiveWorkbook

# Helper: write a string value into a cell while keeping it text (so things
# like leading-zero fund codes "000906" and decimal-looking strings "20.45"
# are preserved verbatim instead of being coerced into numbers).
function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right before the existing "2022-Q3"
#    worksheet (i.e. right after "总计").
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# Header row (matches every other quarterly detail sheet).
$q4Sheet.Cells.Item(1, 2).Value = "基金代码"
$q4Sheet.Cells.Item(1, 3).Value = "基金名称"
$q4Sheet.Cells.Item(1, 4).Value = "基金规模"
$q4Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q4Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1, 8).Value = "仓位排名"

# Row 2: 270023 / 广发全球精选股票（QDII）
$q4Sheet.Cells.Item(2, 1).Value = 0
Set-TextCell $q4Sheet 2 2 "270023"
Set-TextCell $q4Sheet 2 3 "广发全球精选股票（QDII）"
Set-TextCell $q4Sheet 2 4 "20.45"
Set-TextCell $q4Sheet 2 5 "82.63"
Set-TextCell $q4Sheet 2 6 "4.99"
Set-TextCell $q4Sheet 2 7 "1.0205"
$q4Sheet.Cells.Item(2, 8).Value = 9

# Row 3: 000906 / 广发全球精选股票（QDII）美元现汇
$q4Sheet.Cells.Item(3, 1).Value = 1
Set-TextCell $q4Sheet 3 2 "000906"
Set-TextCell $q4Sheet 3 3 "广发全球精选股票（QDII）美元现汇"
Set-TextCell $q4Sheet 3 4 "20.45"
Set-TextCell $q4Sheet 3 5 "82.63"
Set-TextCell $q4Sheet 3 6 "4.99"
Set-TextCell $q4Sheet 3 7 "1.0205"
$q4Sheet.Cells.Item(3, 8).Value = 9

# ---------------------------------------------------------------------------
# 2. Insert a new summary row for "2022-Q4" at the top of the data in the
#    "总计" sheet, pushing every existing row down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 2.04

# Renumber column A (the 0-based row index) for the rows that shifted down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(8, 1).Value = 6

Write-Output "2022-Q4 sheet added; 总计 summary row inserted."
